$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("episodios")
$ws2 = $wb.Worksheets.Item("temporadas")

# --- Update row 3 (episodio-1_temporada-1): hora_estreno / minuto_estreno ---
$ws.Range("M3").Value = 18
$ws.Range("N3").Value = 0

# --- Update row 4 (episodio-2_temporada-1): fecha_estreno / hora_estreno / minuto_estreno ---
# Build the replacement date text via a formula in a scratch cell, then copy/paste the
# computed value back in, so Excel stores it as plain text instead of auto-parsing it
# into a date serial (matches the original "dd/mmm/yyyy"-as-text convention used in this sheet).
$scratch = $ws.Range("ZZ1")
$scratch.Formula = '=CONCATENATE("23","/jun/","2023")'
$scratch.Copy()
$ws.Range("L4").PasteSpecial(-4163)
$scratch.ClearContents()

$ws.Range("M4").Value = 18
$ws.Range("N4").Value = 30

# --- Restore cursor / selection state ---
$ws.Activate()
$ws.Range("L4").Select()

$ws2.Activate()
$ws2.Range("A2").Select()

$ws.Activate()
